$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4977.6665
$ws.Range("I74").Value = 4973.75
$ws.Range("K74").Value = 4973.75
$ws.Range("M74").Value = -4037.75
$ws.Range("H77").Value = 4977.6665
$ws.Range("I77").Value = 4973.75
$ws.Range("K77").Value = 24868.75
$ws.Range("M77").Value = -20188.75
$ws.Range("H92").Value = 790.1
$ws.Range("I92").Value = 783
$ws.Range("K92").Value = 783
$ws.Range("M92").Value = 465
$ws.Range("H116").Value = 3753
$ws.Range("J116").Value = 3006
$ws.Range("L116").Value = 3006
$ws.Range("N116").Value = -9890
$ws.Range("H132").Value = 4103691
$ws.Range("I132").Value = 4550748.5
$ws.Range("K132").Value = 13652245.5
$ws.Range("M132").Value = -13649715.5
$ws.Range("H137").Value = 1392.6757
$ws.Range("I137").Value = 1105.56
$ws.Range("J137").Value = 1990.8334
$ws.Range("K137").Value = 3316.68
$ws.Range("L137").Value = 5972.5002
$ws.Range("M137").Value = -766.6799999999998
$ws.Range("N137").Value = -11072.5002
$ws.Range("H138").Value = 2361.01
$ws.Range("I138").Value = 1132.409
$ws.Range("J138").Value = 2707.5386
$ws.Range("K138").Value = 3397.227
$ws.Range("L138").Value = 8122.6158
$ws.Range("M138").Value = 1742.773
$ws.Range("N138").Value = -18402.6158
$ws.Range("H141").Value = 3280.5293
$ws.Range("I141").Value = 2761.818
$ws.Range("J141").Value = 4231.5
$ws.Range("K141").Value = 8285.454000000002
$ws.Range("L141").Value = 12694.5
$ws.Range("M141").Value = -3105.454000000002
$ws.Range("N141").Value = -23054.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24727.9
$ws.Range("I32").Value = 8104.1265
$ws.Range("J32").Value = 134167.75
$ws.Range("K32").Value = 8104.1265
$ws.Range("L32").Value = 134167.75
$ws.Range("M32").Value = -7817.1265
$ws.Range("N32").Value = -134741.75
$ws.Range("H61").Value = 1669.1923
$ws.Range("I61").Value = 1242.9375
$ws.Range("J61").Value = 2351.2
$ws.Range("K61").Value = 1242.9375
$ws.Range("L61").Value = 2351.2
$ws.Range("M61").Value = -1030.9375
$ws.Range("N61").Value = -2775.2
$ws.Range("H97").Value = 44935.695
$ws.Range("I97").Value = 84340.664
$ws.Range("K97").Value = 84340.664
$ws.Range("M97").Value = -83844.664
$ws.Range("H122").Value = 1977.2572
$ws.Range("I122").Value = 1653.4
$ws.Range("J122").Value = 2786.9
$ws.Range("K122").Value = 4960.200000000001
$ws.Range("L122").Value = 8360.700000000001
$ws.Range("M122").Value = -2510.200000000001
$ws.Range("N122").Value = -13260.7
$ws.Range("H132").Value = 24374.893
$ws.Range("I132").Value = 36070
$ws.Range("K132").Value = 108210
$ws.Range("M132").Value = -105680
$ws.Range("H136").Value = 1669.1923
$ws.Range("I136").Value = 1242.9375
$ws.Range("J136").Value = 2351.2
$ws.Range("K136").Value = 3728.8125
$ws.Range("L136").Value = 7053.599999999999
$ws.Range("M136").Value = -1178.8125
$ws.Range("N136").Value = -12153.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 468.26666
$ws.Range("I94").Value = 409.53845
$ws.Range("J94").Value = 850
$ws.Range("K94").Value = 409.53845
$ws.Range("L94").Value = 850
$ws.Range("M94").Value = 41.46154999999999
$ws.Range("N94").Value = -1752
$ws.Range("H99").Value = 1835.3846
$ws.Range("I99").Value = 1286.6666
$ws.Range("J99").Value = 3070
$ws.Range("K99").Value = 1286.6666
$ws.Range("L99").Value = 3070
$ws.Range("M99").Value = 211.3334
$ws.Range("N99").Value = -6066
$ws.Range("H134").Value = 2810.2654
$ws.Range("I134").Value = 2680.7144
$ws.Range("J134").Value = 3587.5715
$ws.Range("K134").Value = 8042.1432
$ws.Range("L134").Value = 10762.7145
$ws.Range("M134").Value = -5507.1432
$ws.Range("N134").Value = -15832.7145
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 40369.95
$ws.Range("J64").Value = 40369.95
$ws.Range("L64").Value = 40369.95
$ws.Range("N64").Value = -40865.95
$ws.Range("H67").Value = 40369.95
$ws.Range("J67").Value = 40369.95
$ws.Range("L67").Value = 40369.95
$ws.Range("N67").Value = -42085.95
$ws.Range("H122").Value = 899.1667
$ws.Range("I122").Value = 863
$ws.Range("J122").Value = 971.5
$ws.Range("K122").Value = 2589
$ws.Range("L122").Value = 2914.5
$ws.Range("M122").Value = -139
$ws.Range("N122").Value = -7814.5
$ws.Range("H132").Value = 2539.5757
$ws.Range("I132").Value = 2457.1724
$ws.Range("J132").Value = 3137
$ws.Range("K132").Value = 7371.5172
$ws.Range("L132").Value = 9411
$ws.Range("M132").Value = -4841.5172
$ws.Range("N132").Value = -14471
$ws.Range("H134").Value = 1407.95
$ws.Range("I134").Value = 697.06665
$ws.Range("J134").Value = 3540.6
$ws.Range("K134").Value = 2091.19995
$ws.Range("L134").Value = 10621.8
$ws.Range("M134").Value = 443.8000499999998
$ws.Range("N134").Value = -15691.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 1132449
$ws.Range("J37").Value = 1132449
$ws.Range("L37").Value = 3397347
$ws.Range("N37").Value = -3397571
$ws.Range("H113").Value = 1008.0417
$ws.Range("J113").Value = 677.0833
$ws.Range("L113").Value = 2031.2499
$ws.Range("N113").Value = -6371.2499
$ws.Range("H131").Value = 517649.75
$ws.Range("J131").Value = 633965.8
$ws.Range("L131").Value = 1901897.4
$ws.Range("N131").Value = -1911977.4
$ws.Range("H132").Value = 372141.6
$ws.Range("I132").Value = 993.38464
$ws.Range("J132").Value = 716779.2
$ws.Range("K132").Value = 8940.46176
$ws.Range("L132").Value = 6451012.8
$ws.Range("M132").Value = -6410.46176
$ws.Range("N132").Value = -6456072.8
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 2062.7812
$ws.Range("I43").Value = 998.08
$ws.Range("K43").Value = 998.08
$ws.Range("M43").Value = -847.08
$ws.Range("H57").Value = 2550
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").Value = $null
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 22389.5
$ws.Range("J5").Value = 22389.5
$ws.Range("L5").Value = 22389.5
$ws.Range("N5").Value = -22615.5
$ws.Range("H20").Value = 28827.428
$ws.Range("J20").Value = 28827.428
$ws.Range("L20").Value = 28827.428
$ws.Range("N20").Value = -29279.428
$ws.Range("H70").Value = 25000
$ws.Range("J70").Value = 25000
$ws.Range("L70").Value = 25000
$ws.Range("N70").Value = -25540
$ws.Range("H73").Value = 25000
$ws.Range("J73").Value = 25000
$ws.Range("L73").Value = 25000
$ws.Range("N73").Value = -26872
$ws.Range("H132").Value = 4390.3335
$ws.Range("J132").Value = 3170.4285
$ws.Range("L132").Value = 9511.2855
$ws.Range("N132").Value = -14571.2855
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").Value = $null
$ws.Range("H28").Value = 9514.286
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 9514.286
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 9514.286
$ws.Range("M28").Value = $null
$ws.Range("N28").Value = -10210.286
$ws.Range("H30").Value = 60009
$ws.Range("I30").Value = 60009
$ws.Range("K30").Value = 60009
$ws.Range("M30").Value = -59902
$ws.Range("H122").Value = 1927.1111
$ws.Range("I122").Value = 1927.1111
$ws.Range("K122").Value = 5781.3333
$ws.Range("M122").Value = -3331.3333
$ws.Range("H136").Value = 18589.139
$ws.Range("I136").Value = 36468.215
$ws.Range("K136").Value = 109404.645
$ws.Range("M136").Value = -106854.645
